$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 2877.4285
$ws.Range("I5").Value = 2877.4285
$ws.Range("K5").Value = 2877.4285
$ws.Range("M5").Value = -2762.4285
# Row 26
$ws.Range("H26").Value = 34392
$ws.Range("I26").Value = 1199.5
$ws.Range("J26").Value = 100777
$ws.Range("K26").Value = 1199.5
$ws.Range("L26").Value = 100777
$ws.Range("M26").Value = -855.5
$ws.Range("N26").Value = -101465
# Row 39
$ws.Range("H39").Value = 210.41667
$ws.Range("I39").Value = 228.72728
$ws.Range("J39").Value = 9
$ws.Range("K39").Value = 686.18184
$ws.Range("L39").Value = 27
$ws.Range("M39").Value = -390.18184
$ws.Range("N39").Value = -619
# Row 52
$ws.Range("H52").Value = 907.1429000000001
$ws.Range("I52").Value = 500
$ws.Range("J52").Value = 975
$ws.Range("K52").Value = 1500
$ws.Range("L52").Value = 2925
$ws.Range("M52").Value = -1340
$ws.Range("N52").Value = -3245
# Row 87
$ws.Range("H87").Value = 49999.332
$ws.Range("J87").Value = 49999.332
$ws.Range("L87").Value = 49999.332
$ws.Range("N87").Value = -52495.332
# Row 88
$ws.Range("H88").Value = 7674.5
$ws.Range("J88").Value = 8442.429
$ws.Range("L88").Value = 8442.429
$ws.Range("N88").Value = -9254.429
# Row 90
$ws.Range("H90").Value = 49999.332
$ws.Range("J90").Value = 49999.332
$ws.Range("L90").Value = 149997.996
$ws.Range("N90").Value = -162477.996
# Row 91
$ws.Range("H91").Value = 7674.5
$ws.Range("J91").Value = 8442.429
$ws.Range("L91").Value = 8442.429
$ws.Range("N91").Value = -11250.429
# Row 98
$ws.Range("H98").Value = 4998.3335
$ws.Range("I98").Value = 4998.5
$ws.Range("J98").Value = 4998
$ws.Range("K98").Value = 4998.5
$ws.Range("L98").Value = 4998
$ws.Range("M98").Value = -3500.5
$ws.Range("N98").Value = -7994
# Row 107
$ws.Range("H107").Value = 2033.6471
$ws.Range("I107").Value = 989.7692
$ws.Range("J107").Value = 5426.25
$ws.Range("K107").Value = 989.7692
$ws.Range("L107").Value = 5426.25
$ws.Range("M107").Value = 930.2308
$ws.Range("N107").Value = -9266.25
# Row 115
$ws.Range("H115").Value = 330
$ws.Range("I115").Value = 330
$ws.Range("K115").Value = 990
$ws.Range("M115").Value = 577
# Row 122
$ws.Range("H122").Value = 4998.3335
$ws.Range("I122").Value = 4998.5
$ws.Range("J122").Value = 4998
$ws.Range("K122").Value = 14995.5
$ws.Range("L122").Value = 14994
$ws.Range("M122").Value = -12545.5
$ws.Range("N122").Value = -19894
# Row 127
$ws.Range("H127").Value = 4737.375
$ws.Range("I127").Value = 3224.75
$ws.Range("J127").Value = 6250
$ws.Range("K127").Value = 9674.25
$ws.Range("L127").Value = 18750
$ws.Range("M127").Value = -4714.25
$ws.Range("N127").Value = -28670
# Row 129
$ws.Range("H129").Value = 3375.6875
$ws.Range("I129").Value = 1259.1666
$ws.Range("J129").Value = 4645.6
$ws.Range("K129").Value = 3777.4998
$ws.Range("L129").Value = 13936.8
$ws.Range("M129").Value = 1222.5002
$ws.Range("N129").Value = -23936.8
# Row 137
$ws.Range("H137").Value = 2670.1667
$ws.Range("J137").Value = 3463.6667
$ws.Range("L137").Value = 10391.0001
$ws.Range("N137").Value = -15491.0001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1112.9166
$ws.Range("I2").Value = 955.5
$ws.Range("J2").Value = 1900
$ws.Range("K2").Value = 955.5
$ws.Range("L2").Value = 1900
$ws.Range("M2").Value = -842.5
$ws.Range("N2").Value = -2126
# Row 53
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("N53").ClearContents()
# Row 88
$ws.Range("H88").Value = 7247.5
$ws.Range("J88").Value = 7570
$ws.Range("L88").Value = 7570
$ws.Range("N88").Value = -8382
# Row 91
$ws.Range("H91").Value = 7247.5
$ws.Range("J91").Value = 7570
$ws.Range("L91").Value = 7570
$ws.Range("N91").Value = -10378
# Row 105
$ws.Range("H105").Value = 60000
$ws.Range("J105").Value = 60000
$ws.Range("L105").Value = 60000
$ws.Range("N105").Value = -66988
# Row 116
$ws.Range("H116").Value = 1112.9166
$ws.Range("I116").Value = 955.5
$ws.Range("J116").Value = 1900
$ws.Range("K116").Value = 955.5
$ws.Range("L116").Value = 1900
$ws.Range("M116").Value = 1338.5
$ws.Range("N116").Value = -6488

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1112.9166
$ws.Range("I3").Value = 955.5
$ws.Range("J3").Value = 1900
$ws.Range("K3").Value = 955.5
$ws.Range("L3").Value = 1900
$ws.Range("M3").Value = -841.5
$ws.Range("N3").Value = -2128
# Row 20
$ws.Range("H20").Value = 1545.5
$ws.Range("I20").Value = 1365.1875
$ws.Range("K20").Value = 1365.1875
$ws.Range("M20").Value = -1118.1875
# Row 82
$ws.Range("H82").Value = 30857.588
$ws.Range("J82").Value = 39997.5
$ws.Range("L82").Value = 39997.5
$ws.Range("N82").Value = -40763.5
# Row 85
$ws.Range("H85").Value = 30857.588
$ws.Range("J85").Value = 39997.5
$ws.Range("L85").Value = 39997.5
$ws.Range("N85").Value = -42649.5
# Row 99
$ws.Range("H99").Value = 3332.75
$ws.Range("I99").Value = 3544.9092
$ws.Range("K99").Value = 3544.9092
$ws.Range("M99").Value = -2046.9092
# Row 134
$ws.Range("H134").Value = 5636.826
$ws.Range("I134").Value = 5688.5
$ws.Range("K134").Value = 17065.5
$ws.Range("M134").Value = -14530.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 17
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()
# Row 55
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
# Row 62
$ws.Range("H62").Value = 9431.666999999999
$ws.Range("I62").Value = 9418
$ws.Range("K62").Value = 9418
$ws.Range("M62").Value = -8794
# Row 65
$ws.Range("H65").Value = 9431.666999999999
$ws.Range("I65").Value = 9418
$ws.Range("K65").Value = 47090
$ws.Range("M65").Value = -43970
# Row 107
$ws.Range("H107").Value = 378.9
$ws.Range("I107").Value = 176.2
$ws.Range("K107").Value = 176.2
$ws.Range("M107").Value = 1743.8
# Row 122
$ws.Range("H122").Value = 2877.1765
$ws.Range("I122").Value = 3222.2856
$ws.Range("J122").Value = 1266.6666
$ws.Range("K122").Value = 9666.856800000001
$ws.Range("L122").Value = 3799.9998
$ws.Range("M122").Value = -7216.856800000001
$ws.Range("N122").Value = -8699.9998

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 1549.52
$ws.Range("I4").Value = 1139.875
$ws.Range("J4").Value = 2277.7778
$ws.Range("K4").Value = 3419.625
$ws.Range("L4").Value = 6833.3334
$ws.Range("M4").Value = -3307.625
$ws.Range("N4").Value = -7057.3334
# Row 37
$ws.Range("H37").Value = 95000
$ws.Range("J37").Value = 95000
$ws.Range("L37").Value = 285000
$ws.Range("N37").Value = -285224
# Row 93
$ws.Range("H93").Value = 2498
$ws.Range("J93").Value = 2498
$ws.Range("L93").Value = 7494
$ws.Range("N93").Value = -11238
# Row 122
$ws.Range("H122").Value = 1174.0476
$ws.Range("J122").Value = 1223
$ws.Range("L122").Value = 11007
$ws.Range("N122").Value = -15907
# Row 129
$ws.Range("H129").Value = 2365.9092
$ws.Range("J129").Value = 3185.1428
$ws.Range("L129").Value = 9555.428400000001
$ws.Range("N129").Value = -19555.4284
# Row 134
$ws.Range("H134").Value = 976.3333
$ws.Range("I134").Value = 976.3333
$ws.Range("K134").Value = 2928.9999
$ws.Range("M134").Value = 2141.0001
# Row 140
$ws.Range("H140").Value = 1497.5
$ws.Range("I140").Value = 1000
$ws.Range("K140").Value = 3000
$ws.Range("M140").Value = 2180

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 7623.75
$ws.Range("I70").Value = 7623.75
$ws.Range("K70").Value = 7623.75
$ws.Range("M70").Value = -7353.75
# Row 73
$ws.Range("H73").Value = 7623.75
$ws.Range("I73").Value = 7623.75
$ws.Range("K73").Value = 7623.75
$ws.Range("M73").Value = -6687.75
# Row 102
$ws.Range("H102").Value = 1609.75
$ws.Range("I102").Value = 1609.75
$ws.Range("K102").Value = 1609.75
$ws.Range("M102").Value = 12.25
# Row 107
$ws.Range("H107").Value = 229.83333
$ws.Range("I107").Value = 245.8
$ws.Range("K107").Value = 245.8
$ws.Range("M107").Value = 1674.2
# Row 126
$ws.Range("H126").Value = 1670.6154
$ws.Range("I126").Value = 1670.6154
$ws.Range("K126").Value = 5011.8462
$ws.Range("M126").Value = -2541.8462

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 31
$ws.Range("H31").Value = 1102.8
$ws.Range("I31").Value = 657
$ws.Range("K31").Value = 657
$ws.Range("M31").Value = -409
